$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update odds that moved between scrapes (rows 3-16) ---
$ws.Cells.Item(3, 7).Value = 1.75
$ws.Cells.Item(3, 9).Value = 4.75
$ws.Cells.Item(3, 17).Value = 2.06
$ws.Cells.Item(3, 18).Value = 1.84
$ws.Cells.Item(3, 26).Value = 13
$ws.Cells.Item(3, 50).Value = 6.5
$ws.Cells.Item(4, 7).Value = 2.4
$ws.Cells.Item(4, 9).Value = 3.4
$ws.Cells.Item(4, 10).Value = 3.25
$ws.Cells.Item(4, 13).Value = 1.11
$ws.Cells.Item(4, 14).Value = 6.5
$ws.Cells.Item(4, 15).Value = 1.5
$ws.Cells.Item(4, 16).Value = 2.63
$ws.Cells.Item(4, 17).Value = 2.63
$ws.Cells.Item(4, 18).Value = 1.5
$ws.Cells.Item(4, 34).Value = 7.5
$ws.Cells.Item(4, 47).Value = 9
$ws.Cells.Item(5, 19).Value = 1.62
$ws.Cells.Item(6, 7).Value = 2.63
$ws.Cells.Item(6, 9).Value = 2.9
$ws.Cells.Item(6, 10).Value = 3.4
$ws.Cells.Item(6, 12).Value = 3.75
$ws.Cells.Item(6, 19).Value = 1.57
$ws.Cells.Item(6, 21).Value = 2.05
$ws.Cells.Item(6, 22).Value = 1.7
$ws.Cells.Item(6, 26).Value = 26
$ws.Cells.Item(6, 33).Value = 501
$ws.Cells.Item(6, 34).Value = 7.5
$ws.Cells.Item(6, 36).Value = 12
$ws.Cells.Item(6, 38).Value = 29
$ws.Cells.Item(6, 52).Value = 34
$ws.Cells.Item(7, 7).Value = 2.3
$ws.Cells.Item(7, 8).Value = 2.75
$ws.Cells.Item(7, 9).Value = 3.5
$ws.Cells.Item(7, 10).Value = 3.25
$ws.Cells.Item(7, 12).Value = 4.5
$ws.Cells.Item(7, 24).Value = 9.5
$ws.Cells.Item(7, 26).Value = 23
$ws.Cells.Item(7, 27).Value = 29
$ws.Cells.Item(7, 50).Value = 5
$ws.Cells.Item(7, 53).Value = 81
$ws.Cells.Item(9, 17).Value = 1.73
$ws.Cells.Item(9, 18).Value = 2.08
$ws.Cells.Item(10, 37).Value = 51
$ws.Cells.Item(10, 41).Value = 8.5
$ws.Cells.Item(11, 13).Value = 1.05
$ws.Cells.Item(11, 14).Value = 11
$ws.Cells.Item(11, 17).Value = 1.84
$ws.Cells.Item(11, 18).Value = 2.06
$ws.Cells.Item(11, 19).Value = 1.36
$ws.Cells.Item(12, 9).Value = 3.75
$ws.Cells.Item(12, 19).Value = 1.53
$ws.Cells.Item(12, 20).Value = 2.38
$ws.Cells.Item(12, 21).Value = 2.1
$ws.Cells.Item(12, 22).Value = 1.67
$ws.Cells.Item(12, 23).Value = 6
$ws.Cells.Item(12, 32).Value = 67
$ws.Cells.Item(12, 36).Value = 13
$ws.Cells.Item(12, 38).Value = 34
$ws.Cells.Item(12, 39).Value = 41
$ws.Cells.Item(12, 42).Value = 26
$ws.Cells.Item(12, 44).Value = 67
$ws.Cells.Item(12, 46).Value = 2.38
$ws.Cells.Item(12, 47).Value = 9
$ws.Cells.Item(12, 48).Value = 67
$ws.Cells.Item(12, 52).Value = 34
$ws.Cells.Item(13, 14).Value = 13
$ws.Cells.Item(13, 17).Value = 1.7
$ws.Cells.Item(13, 18).Value = 2.1
$ws.Cells.Item(14, 9).Value = 2.8
$ws.Cells.Item(14, 34).Value = 7.5
$ws.Cells.Item(14, 39).Value = 34
$ws.Cells.Item(14, 40).Value = 4.5
$ws.Cells.Item(16, 7).Value = 2.4
$ws.Cells.Item(16, 9).Value = 2.88
$ws.Cells.Item(16, 10).Value = 3.1
$ws.Cells.Item(16, 11).Value = 2.1
$ws.Cells.Item(16, 13).Value = 1.06
$ws.Cells.Item(16, 14).Value = 10
$ws.Cells.Item(16, 25).Value = 10
$ws.Cells.Item(16, 28).Value = 29
$ws.Cells.Item(16, 34).Value = 9.5
$ws.Cells.Item(16, 50).Value = 4.75
$ws.Cells.Item(16, 51).Value = 15

# --- 2) Remove the three stale fixtures (old rows 17-19: Germany 3.Liga, Greece Super League, Romania Liga 1) ---
$ws.Rows("17:19").Delete()

# --- 3) Insert a fresh row for the new fixture ahead of the two remaining fixtures (now at 17/18) ---
$ws.Rows(17).Insert()

# --- 4) Populate new row 17: PARAGUAY - PRIMERA DIVISION (Olimpia Asuncion vs Nacional Asuncion) ---
$ws.Cells.Item(17, 1).Value = "2Zk1gpyH"
$ws.Cells.Item(17, 2).Value = "24/11/2024"
$ws.Cells.Item(17, 3).Value = "18:00"
$ws.Cells.Item(17, 4).Value = "PARAGUAY - PRIMERA DIVISION"
$ws.Cells.Item(17, 5).Value = "Olimpia Asuncion"
$ws.Cells.Item(17, 6).Value = "Nacional Asuncion"
$ws.Cells.Item(17, 7).Value = 1.65
$ws.Cells.Item(17, 8).Value = 3.5
$ws.Cells.Item(17, 9).Value = 5
$ws.Cells.Item(17, 10).Value = 2.38
$ws.Cells.Item(17, 11).Value = 2.05
$ws.Cells.Item(17, 12).Value = 6
$ws.Cells.Item(17, 13).Value = 1.08
$ws.Cells.Item(17, 14).Value = 7.5
$ws.Cells.Item(17, 15).Value = 1.44
$ws.Cells.Item(17, 16).Value = 2.63
$ws.Cells.Item(17, 17).Value = 2.35
$ws.Cells.Item(17, 18).Value = 1.57
$ws.Cells.Item(17, 19).Value = 1.5
$ws.Cells.Item(17, 20).Value = 2.5
$ws.Cells.Item(17, 21).Value = 2.25
$ws.Cells.Item(17, 22).Value = 1.57
$ws.Cells.Item(17, 23).Value = 5.5
$ws.Cells.Item(17, 24).Value = 7
$ws.Cells.Item(17, 25).Value = 9
$ws.Cells.Item(17, 26).Value = 12
$ws.Cells.Item(17, 27).Value = 17
$ws.Cells.Item(17, 28).Value = 34
$ws.Cells.Item(17, 29).Value = 7.5
$ws.Cells.Item(17, 30).Value = 7
$ws.Cells.Item(17, 31).Value = 21
$ws.Cells.Item(17, 32).Value = 81
$ws.Cells.Item(17, 33).Value = 351
$ws.Cells.Item(17, 34).Value = 11
$ws.Cells.Item(17, 35).Value = 23
$ws.Cells.Item(17, 36).Value = 17
$ws.Cells.Item(17, 37).Value = 51
$ws.Cells.Item(17, 38).Value = 41
$ws.Cells.Item(17, 39).Value = 51
$ws.Cells.Item(17, 40).Value = 3.5
$ws.Cells.Item(17, 41).Value = 9
$ws.Cells.Item(17, 42).Value = 23
$ws.Cells.Item(17, 43).Value = 34
$ws.Cells.Item(17, 44).Value = 67
$ws.Cells.Item(17, 45).Value = 251
$ws.Cells.Item(17, 46).Value = 2.5
$ws.Cells.Item(17, 47).Value = 9.5
$ws.Cells.Item(17, 48).Value = 81
$ws.Cells.Item(17, 50).Value = 7
$ws.Cells.Item(17, 51).Value = 34
$ws.Cells.Item(17, 52).Value = 41
$ws.Cells.Item(17, 53).Value = 126
$ws.Cells.Item(17, 54).Value = 151
$ws.Cells.Item(17, 55).Value = 351

# --- 5) Row 18 (SPAIN - LALIGA2, Elche vs R. Oviedo) keeps its identity; only these odds moved ---
$ws.Cells.Item(18, 7).Value = 2.05
$ws.Cells.Item(18, 9).Value = 4
$ws.Cells.Item(18, 10).Value = 2.88
$ws.Cells.Item(18, 12).Value = 4.75
$ws.Cells.Item(18, 13).Value = 1.11
$ws.Cells.Item(18, 14).Value = 6.5
$ws.Cells.Item(18, 15).Value = 1.5
$ws.Cells.Item(18, 16).Value = 2.5
$ws.Cells.Item(18, 24).Value = 8.5
$ws.Cells.Item(18, 27).Value = 21
$ws.Cells.Item(18, 34).Value = 8.5
$ws.Cells.Item(18, 44).Value = 81
$ws.Cells.Item(18, 50).Value = 5.5
$ws.Cells.Item(18, 51).Value = 23
$ws.Cells.Item(18, 53).Value = 81
$ws.Cells.Item(18, 54).Value = 126

# Row 19 (USA - MLS, Orlando City vs Atlanta Utd) is unchanged after the shift; nothing further to do.
